$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their original text (string) representation,
# matching the source data which stores numeric-looking price/percent
# values as plain text (inline strings), not as numbers.
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.937.74"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.895.96"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "0.7753"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("D6").Value = "244.78"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.3139"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "25.78"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "0.08939"
$ws.Range("E11").Value = "  +10.28%  "
$ws.Range("D12").Value = "0.7739"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "5.457"
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("D14").Value = "94.80"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").Value = "1.815.51"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").Value = "6.200"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "29.938.42"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "13.99"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "246.92"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "0.000007890"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").Value = "8.157"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "2.133.76"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "0.1589"
$ws.Range("E25").Value = "  -5.36%  "
$ws.Range("D26").Value = "9.549"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").Value = "163.22"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").Value = "18.85"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").Value = "2.050"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").Value = "1.427"
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("D31").Value = "1.545"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "4.548"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "4.124"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "0.05527"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "1.250"
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("D36").Value = "0.7548"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").Value = "0.9976"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "2.719"
$ws.Range("E38").Value = "  +3.37%  "
$ws.Range("D39").Value = "0.01968"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "0.4522"
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("D42").Value = "73.95"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").Value = "6.064"
$ws.Range("E43").Value = "  +2.52%  "
$ws.Range("D44").Value = "1.086.31"
$ws.Range("E44").Value = "  -5.97%  "
$ws.Range("D45").Value = "0.8555"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D47").Value = "1.898"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "102.86"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").Value = "7.620"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").Value = "9.874"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "2.997"
$ws.Range("E51").Value = "  -1.51%  "
